# Update the cryptocurrency price/volume table (and reorder a few coin rows)
# to match the refreshed data pulled by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.982.92'
$ws.Range("E2").Value = '  +4.21%  '
$ws.Range("D3").Value = '2.240.63'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.96'
$ws.Range("E5").Value = '  +3.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.614'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.41'
$ws.Range("E7").Value = '  +8.70%  '
$ws.Range("E9").Value = '  +6.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.08'
$ws.Range("E10").Value = '  +1.88%  '
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.53'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.96'
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("E14").Value = '  +0.87%  '
$ws.Range("D15").Value = '2.556.77'
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.67'
$ws.Range("E16").Value = '  +5.38%  '
$ws.Range("D17").Value = '2.245.52'
$ws.Range("E17").Value = '  +3.60%  '
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("D19").Value = '42.912.84'
$ws.Range("E19").Value = '  +4.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000104'
$ws.Range("E20").Value = '  +2.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.34'
$ws.Range("E21").Value = '  +1.11%  '
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.22'
$ws.Range("E23").Value = '  +4.51%  '
$ws.Range("E24").Value = '  +14.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '230.92'
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.88'
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.37'
$ws.Range("E28").Value = '  -4.95%  '
$ws.Range("E29").Value = '  +2.27%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.22'
$ws.Range("E31").Value = '  +24.56%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.75'
$ws.Range("E32").Value = '  +3.31%  '
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0796'
$ws.Range("E34").Value = '  +3.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.37'
$ws.Range("E35").Value = '  +3.92%  '
$ws.Range("E36").Value = '  +11.47%  '
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.39'
$ws.Range("E38").Value = '  +6.12%  '
$ws.Range("E39").Value = '  +14.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.97'
$ws.Range("E40").Value = '  +8.47%  '
$ws.Range("E41").Value = '  +3.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.58'
$ws.Range("E42").Value = '  +2.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.205'
$ws.Range("E43").Value = '  +7.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.12'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '106.16'
$ws.Range("E45").Value = '  +8.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.68'
$ws.Range("E46").Value = '  +4.24%  '
$ws.Range("E47").Value = '  +1.31%  '
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("B49").Value = 'WOONetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.441'
$ws.Range("E49").Value = '  +15.43%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.31'
$ws.Range("E50").Value = '  +3.98%  '
$ws.Range("E51").Value = '  +1.01%  '
